$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grouped matches")

$data = @(
    @(2, 'P2614', '{''eft:nyen-lotsawa-darma-drak''}'),
    @(3, 'P8261', '{''eft:munivarman'', ''eft:munivarma''}'),
    @(4, 'P00KG07267', '{''eft:sarvajnadeva'', ''eft:sarvanyadeva''}'),
    @(5, 'P2551', '{''eft:blo-ldan-shes-rab''}'),
    @(6, 'P8273', '{''eft:rin-chen-tsho'', ''eft:rinchen-tso''}'),
    @(7, 'P5788', '{''eft:th-nmi-sambhota''}'),
    @(8, 'P8213', '{''eft:t-vidyakarasimha'', ''eft:vidyakarasimha''}'),
    @(9, 'P2548', '{''eft:prajnavarman'', ''eft:prajnavarma''}'),
    @(10, 'P8222', '{''eft:jnanasidhi'', ''eft:jnanasiddhi''}'),
    @(11, 'P8268', '{''eft:buddhaprabha''}'),
    @(12, 'P8182', '{''eft:band-paltsek'', ''eft:kawa-paltsek-under-the-name-paltsek-raksita-'', ''eft:ska-ba-dpal-brtsegs'', ''eft:dpal-brtsegs'', ''eft:paltsek'', ''eft:ban-de-dpal-brtsegs''}'),
    @(13, 'P3709', '{''eft:phakpa-sherab''}'),
    @(14, 'P3214', '{''eft:danasila''}'),
    @(15, 'P4263', '{''eft:dge-ba-dpal'', ''eft:gew-pal''}'),
    @(16, 'P3456', '{''eft:tshul-khrims-rgyal-ba'', ''eft:tsultrim-gyalwa''}'),
    @(17, 'P8217', '{''eft:jnanagarbha'', ''eft:t-jnanagarbha''}'),
    @(18, 'P8220', '{''eft:devacandra''}'),
    @(19, 'P8249', '{''eft:pandita-dharmakara'', ''eft:dharmakara''}'),
    @(20, 'P0TMP092', '{''eft:anandasri-s-''}'),
    @(21, 'P8271', '{''eft:kumararaksita''}'),
    @(22, 'P0TMPT007', '{''eft:rnam-par-mi-rtog-pa''}'),
    @(23, 'P4259', '{''eft:palgyi-lh-npo'', ''eft:ban-de-dpal-gyi-lhun-po'', ''eft:dpal-gyi-lhun-po''}'),
    @(24, 'P8263', '{''eft:leki-d-''}'),
    @(25, 'P8276', '{''eft:wang-phab-zhwun-wang-phan-zhun-''}'),
    @(26, 'P8212', '{''eft:devendraraksita''}'),
    @(27, 'P4258', '{''eft:dpal-byor''}'),
    @(28, 'P3285', '{''eft:sakya-yesh-''}'),
    @(29, 'P8265', '{''eft:ratnaraksita''}'),
    @(30, 'P0TMP080', '{''eft:hwa-shang-zab-mo''}'),
    @(31, 'P0RK8', '{''eft:dharmapala''}'),
    @(32, 'P8216', '{''eft:sakya-lodr-''}'),
    @(33, 'P4256', '{''eft:lotsawa-zangkyong-bzang-skyong-''}'),
    @(34, 'P8277', '{''eft:rgya-mtsho-i-sde''}'),
    @(35, 'P2956', '{''eft:krsnapandita''}'),
    @(36, '?', '{''eft:sherap-'', ''eft:sakyasena'', ''eft:vajrvisramitra''}'),
    @(37, 'P4CZ16819', '{''eft:sakyaprabha''}'),
    @(38, 'P8245', '{''eft:buddhakaravarma''}'),
    @(39, 'P8221', '{''eft:g-ch-drup''}'),
    @(40, 'P3890', '{''eft:ch-kyi-sherab''}'),
    @(41, 'P2637', '{''eft:trakpa-gyaltsen''}'),
    @(42, 'P4CZ16780', '{''eft:manjusrigarbha''}'),
    @(43, 'P8171', '{''eft:dharmasribhadra''}'),
    @(44, 'P8211', '{''eft:bidyakaraprabha'', ''eft:vidyakaraprabha''}'),
    @(45, 'P8260', '{''eft:dpal-dbyangs''}'),
    @(46, 'P4CZ15137', '{''eft:kumarakalasa''}'),
    @(47, 'P1KG8854', '{''eft:silendrabodhi'', ''eft:srilendrabodhi'', ''eft:surendrabodhi''}'),
    @(48, 'P0TMP098', '{''eft:jinavara''}'),
    @(49, 'P8219', '{''eft:visuddhasimha''}'),
    @(50, 'P8266', '{''eft:lotsawa-band-dharmatasila'', ''eft:dharmatasila'', ''eft:ch-nyi-tsultrim''}'),
    @(51, 'P8205', '{''eft:band-yesh-de'', ''eft:ye-shes-sde'', ''eft:yesh-de'', ''eft:zhang-yesh-d-'', ''eft:band-yesh-d-'', ''eft:yesh-d-ye-shes-sde-'', ''eft:yesh-d-''}'),
    @(52, 'P4255', '{''eft:t-jnanagarbha'', ''eft:yesh-nyingpo'', ''eft:ye-shes-snying-po''}'),
    @(53, 'P1321', '{''eft:shang-buchikpa''}'),
    @(54, 'P8209', '{''eft:dzi-na-mi-tra-k-'', ''eft:jinamitra'', ''eft:jinamitra-k-''}'),
    @(55, 'P0TMP104', '{''eft:punyasambhava''}'),
    @(56, 'P8267', '{''eft:vijayasila''}'),
    @(57, 'P4CZ15308', '{''eft:vairocanaraksita''}'),
    @(58, 'P3379', '{''eft:dipamkarasrijnana'', ''eft:dipamkara-srijnana''}'),
    @(59, 'P8183', '{''eft:klu-i-rgyal-mtshan'', ''eft:cog-ro-klu-i-rgyal-mtshan''}'),
    @(60, 'P753', '{''eft:rin-chen-bzag-po'', ''eft:rin-chen-bzang-po'', ''eft:rinchen-zangpo''}'),
    @(61, 'P8093', '{''eft:kamalagupta''}'),
    @(62, 'P8206', '{''eft:celu''}'),
    @(63, 'P6453', '{''eft:tsultrim-gyaltsen''}'),
    @(64, 'https://lod.dila.edu.tw/resource.php?id=A000089', '{''eft:siladharma''}'),
    @(65, 'P8151', '{''eft:gayadhara''}'),
    @(66, 'P4242', '{''eft:sherab-lekpa''}'),
    @(67, 'P5651', '{''eft:pa-tshab-nyi-ma-grags'', ''eft:patsap-nyima-drak-''}'),
    @(68, 'P8228', '{''eft:surendrabodhi''}'),
    @(69, 'P8269', '{''eft:dgon-gling-rma''}'),
    @(70, 'P8278', '{''eft:dge-ba-i-blo-gros'', ''eft:gewai-lodr-''}'),
    @(71, 'P2557', '{''eft:-brom'', ''eft:-brom-rgyal-ba-i-byung-gnas''}'),
    @(72, 'P1242', '{''eft:g-wai-lodr-''}'),
    @(73, 'P8280', '{''eft:subhasita''}'),
    @(74, 'P3458', '{''eft:g-lhets-''}'),
)

foreach ($row in $data) {
    $r = $row[0]
    $b = $row[1]
    $c = $row[2]
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
}
